$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header columns and fix title-casing of Spanish connector words ---
# (de/del/la/las/el/los/y -> De/Del/La/Las/El/Los/Y) in municipality/state names,
# plus a couple of one-off case fixes (GUANAJUATO -> Guanajuato, MonteMorelos -> Montemorelos)
$ws.Range("A1").Value = 'mx_state'
$ws.Range("B1").Value = 'mx_municipality'
$ws.Range("C1").Value = 'n_matriculas'
$ws.Range("D1").Value = 'pct_matriculas'
$ws.Range("B7").Value = 'Pabellón De Arteaga'
$ws.Range("B8").Value = 'Rincón De Romos'
$ws.Range("B13").Value = 'Playas De Rosarito'
$ws.Range("B31").Value = 'Amatenango De La Frontera'
$ws.Range("B33").Value = 'Bejucal De Ocampo'
$ws.Range("B40").Value = 'Chiapa De Corzo'
$ws.Range("B45").Value = 'Comitán De Domínguez'
$ws.Range("B59").Value = 'Mazapa De Madero'
$ws.Range("B70").Value = 'Salto De Agua'
$ws.Range("B94").Value = 'Hidalgo Del Parral'
$ws.Range("B99").Value = 'San Francisco Del Oro'
$ws.Range("B114").Value = 'San Juan De Sabinas'
$ws.Range("A127").Value = 'Ciudad De México'
$ws.Range("B131").Value = 'Cuajimalpa De Morelos'
$ws.Range("B144").Value = 'Coneto De Comonfort'
$ws.Range("B156").Value = 'Nombre De Dios'
$ws.Range("B163").Value = 'San Juan Del Río'
$ws.Range("A170").Value = 'Estado De México'
$ws.Range("B170").Value = 'Acambay De Ruíz Castañeda'
$ws.Range("B173").Value = 'Almoloya De Alquisiras'
$ws.Range("B174").Value = 'Almoloya De Juárez'
$ws.Range("B181").Value = 'Atizapán De Zaragoza'
$ws.Range("B186").Value = 'Chapa De Mota'
$ws.Range("B189").Value = 'Coacalco De Berriozábal'
$ws.Range("B194").Value = 'Ecatepec De Morelos'
$ws.Range("B196").Value = 'Ixtapan De La Sal'
$ws.Range("B208").Value = 'Naucalpan De Juárez'
$ws.Range("B215").Value = 'San Felipe Del Progreso'
$ws.Range("B217").Value = 'San Simón De Guerero'
$ws.Range("B232").Value = 'Tlalnepantla De Baz'
$ws.Range("B238").Value = 'Valle De Bravo'
$ws.Range("B239").Value = 'Valle De Chalco Solidaridad'
$ws.Range("B240").Value = 'Villa De Allende'
$ws.Range("B241").Value = 'Villa Del Carbón'
$ws.Range("A250").Value = 'Guanajuato'
$ws.Range("B253").Value = 'Apaseo El Alto'
$ws.Range("B254").Value = 'Apaseo El Grande'
$ws.Range("B260").Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range("B264").Value = 'Jaral Del Progreso'
$ws.Range("B272").Value = 'Purísima Del Rincón'
$ws.Range("B275").Value = 'San Diego De La Unión'
$ws.Range("B277").Value = 'San Francisco Del Rincón'
$ws.Range("B279").Value = 'San Luis De La Paz'
$ws.Range("B280").Value = 'Santa Cruz De Juventino Rosas'
$ws.Range("B282").Value = 'Silao De La Victoria'
$ws.Range("B286").Value = 'Valle De Santiago'
$ws.Range("B290").Value = 'Acapulco De Juárez'
$ws.Range("B292").Value = 'Ajuchitlán Del Progreso'
$ws.Range("B293").Value = 'Alcozauca De Guerero'
$ws.Range("B297").Value = 'Atenango Del Río'
$ws.Range("B298").Value = 'Atlamajalcingo Del Monte'
$ws.Range("B300").Value = 'Atoyac De Álvarez'
$ws.Range("B301").Value = 'Ayutla De Los Libres'
$ws.Range("B304").Value = 'Chilapa De Álvarez'
$ws.Range("B305").Value = 'Chilpancingo De Los Bravo'
$ws.Range("B309").Value = 'Coyuca De Benítez'
$ws.Range("B310").Value = 'Coyuca De Catalán'
$ws.Range("B313").Value = 'Cutzamala De Pinzón'
$ws.Range("B319").Value = 'Huitzuco De Los Figueroa'
$ws.Range("B320").Value = 'Iguala De La Independencia'
$ws.Range("B322").Value = 'Zihuatanejo De Azueta'
$ws.Range("B324").Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range("B327").Value = 'Mártir De Cuilapan'
$ws.Range("B339").Value = 'Taxco De Alarcón'
$ws.Range("B341").Value = 'Técpan De Galeana'
$ws.Range("B343").Value = 'Tepecoacuilco De Trujano'
$ws.Range("B345").Value = 'Tixtla De Guerero'
$ws.Range("B347").Value = 'Tlapa De Comonfort'
$ws.Range("B356").Value = 'Agua Blanca De Iturbide'
$ws.Range("B361").Value = 'Atotonilco El Grande'
$ws.Range("B365").Value = 'Cuautepec De Hinojosa'
$ws.Range("B368").Value = 'Huasca De Ocampo'
$ws.Range("B370").Value = 'Huejutla De Reyes'
$ws.Range("B373").Value = 'Jacala De Ledezma'
$ws.Range("B378").Value = 'Mineral Del Chico'
$ws.Range("B379").Value = 'Mineral Del Monte'
$ws.Range("B380").Value = 'Molango De Escamilla'
$ws.Range("B382").Value = 'Omitlán De Juárez'
$ws.Range("B383").Value = 'Pachuca De Soto'
$ws.Range("B386").Value = 'Progreso De Obregón'
$ws.Range("B387").Value = 'Santiago De Anaya'
$ws.Range("B390").Value = 'Tenango De Doria'
$ws.Range("B392").Value = 'Tepehuacán De Guerero'
$ws.Range("B393").Value = 'Tepeji Del Río De Ocampo'
$ws.Range("B398").Value = 'Tula De Allende'
$ws.Range("B399").Value = 'Tulancingo De Bravo'
$ws.Range("B401").Value = 'Zacualtipán De Ángeles'
$ws.Range("B402").Value = 'Zapotlán De Juárez'
$ws.Range("B405").Value = 'Ahualulco De Mercado'
$ws.Range("B409").Value = 'Atotonilco El Alto'
$ws.Range("B410").Value = 'Autlán De Navarro'
$ws.Range("B413").Value = 'Cañadas De Obregón'
$ws.Range("B420").Value = 'Concepción De Buenos Aires'
$ws.Range("B427").Value = 'Encarnación De Díaz'
$ws.Range("B431").Value = 'Huejuquilla El Alto'
$ws.Range("B432").Value = 'Ixtlahuacán Del Río'
$ws.Range("B435").Value = 'Jilotlán De Los Dolores'
$ws.Range("B440").Value = 'La Manzanilla De La Paz'
$ws.Range("B441").Value = 'Lagos De Moreno'
$ws.Range("B448").Value = 'Ojuelos De Jalisco'
$ws.Range("B453").Value = 'San Juan De Los Lagos'
$ws.Range("B455").Value = 'San Miguel El Alto'
$ws.Range("B458").Value = 'Talpa De Allende'
$ws.Range("B460").Value = 'Techaluta De Montenegro'
$ws.Range("B462").Value = 'Tepatitlán De Morelos'
$ws.Range("B464").Value = 'Tizapán El Alto'
$ws.Range("B465").Value = 'Tlajomulco De Zúñiga'
$ws.Range("B472").Value = 'Unión De Tula'
$ws.Range("B473").Value = 'Valle De Guadalupe'
$ws.Range("B474").Value = 'Valle De Juárez'
$ws.Range("B479").Value = 'Yahualica De González Gallo'
$ws.Range("B480").Value = 'Zacoalco De Torres'
$ws.Range("B483").Value = 'Zapotlán Del Rey'
$ws.Range("B484").Value = 'Zapotlán El Grande'
$ws.Range("B564").Value = 'Tiquicheo De Nicolás Romero'
$ws.Range("B588").Value = 'Coatlán Del Río'
$ws.Range("B594").Value = 'Jonacatepec De Leandro Valle'
$ws.Range("B597").Value = 'Puente De Ixtla'
$ws.Range("B602").Value = 'Tlaltizapán De Zapata'
$ws.Range("B613").Value = 'Ixtlán Del Río'
$ws.Range("B618").Value = 'Santa María Del Oro'
$ws.Range("B631").Value = 'Montemorelos'
$ws.Range("B634").Value = 'San Nicolás De Los Garza'
$ws.Range("B637").Value = 'Acatlán De Pérez Figueroa'
$ws.Range("B641").Value = 'Chalcatongo De Hidalgo'
$ws.Range("B644").Value = 'Constancia Del Rosario'
$ws.Range("B646").Value = 'El Barrio De La Soledad'
$ws.Range("B647").Value = 'Fresnillo De Trujano'
$ws.Range("B648").Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range("B649").Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range("B650").Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range("B651").Value = 'Huautla De Jiménez'
$ws.Range("B652").Value = 'Ixtlán De Juárez'
$ws.Range("B653").Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range("B657").Value = 'Mariscala De Juárez'
$ws.Range("B658").Value = 'Mártires De Tacubaya'
$ws.Range("B660").Value = 'Mazatlán Villa De Flores'
$ws.Range("B661").Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range("B662").Value = 'Oaxaca De Juárez'
$ws.Range("B663").Value = 'Ocotlán De Morelos'
$ws.Range("B664").Value = 'Pinotepa De Don Luis'
$ws.Range("B666").Value = 'Putla Villa De Guerero'
$ws.Range("B667").Value = 'Reforma De Pineda'
$ws.Range("B668").Value = 'Rojas De Cuauhtémoc'
$ws.Range("B681").Value = 'San Felipe Jalapa De Díaz'
$ws.Range("B694").Value = 'San Juan Del Estado'
$ws.Range("B733").Value = 'San Pedro Y San Pablo Teposcolula'
$ws.Range("B772").Value = 'Santo Domingo De Morelos'
$ws.Range("B778").Value = 'Teotitlán De Flores Magón'
$ws.Range("B779").Value = 'Tepelmeme Villa De Morelos'
$ws.Range("B780").Value = 'Tlacolula De Matamoros'
$ws.Range("B781").Value = 'Totontepec Villa De Morelos'
$ws.Range("B784").Value = 'Villa De Chilapa De Díaz'
$ws.Range("B785").Value = 'Villa De Etla'
$ws.Range("B786").Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range("B787").Value = 'Villa De Zaachila'
$ws.Range("B788").Value = 'Villa Sola De Vega'
$ws.Range("B789").Value = 'Zimatlán De Álvarez'
$ws.Range("B804").Value = 'Chalchicomula De Sesma'
$ws.Range("B812").Value = 'Chila De La Sal'
$ws.Range("B820").Value = 'Cuayuca De Andrade'
$ws.Range("B821").Value = 'Cuetzalan Del Progreso'
$ws.Range("B835").Value = 'Huehuetlán El Chico'
$ws.Range("B836").Value = 'Huehuetlán El Grande'
$ws.Range("B839").Value = 'Ixcamilpa De Guerero'
$ws.Range("B841").Value = 'Izúcar De Matamoros'
$ws.Range("B849").Value = 'Los Reyes De Juárez'
$ws.Range("B855").Value = 'Palmar De Bravo'
$ws.Range("B872").Value = 'San Salvador El Seco'
$ws.Range("B876").Value = 'Tecali De Herrera'
$ws.Range("B882").Value = 'Tepanco De López'
$ws.Range("B887").Value = 'Tepexi De Rodríguez'
$ws.Range("B889").Value = 'Tepeyahualco De Cuauhtémoc'
$ws.Range("B890").Value = 'Tetela De Ocampo'
$ws.Range("B895").Value = 'Tlacotepec De Benito Juárez'
$ws.Range("B916").Value = 'Amealco De Bonfil'
$ws.Range("B918").Value = 'Cadereyta De Montes'
$ws.Range("B921").Value = 'Jalpan De Serra'
$ws.Range("B922").Value = 'Landa De Matamoros'
$ws.Range("B924").Value = 'Pinal De Amoles'
$ws.Range("B926").Value = 'San Juan Del Río'
$ws.Range("B938").Value = 'Ciudad Del Maíz'
$ws.Range("B945").Value = 'Mexquitic De Carmona'
$ws.Range("B953").Value = 'Santa María Del Río'
$ws.Range("B954").Value = 'Soledad De Graciano Sánchez'
$ws.Range("B958").Value = 'Tanquián De Escobedo'
$ws.Range("B961").Value = 'Villa De Arista'
$ws.Range("B962").Value = 'Villa De Arriaga'
$ws.Range("B963").Value = 'Villa De Ramos'
$ws.Range("B964").Value = 'Villa De Reyes'
$ws.Range("B1001").Value = 'Jalpa De Méndez'
$ws.Range("B1027").Value = 'Soto La Marina'
$ws.Range("B1043").Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range("B1046").Value = 'Papalotla De Xicohténcatl'
$ws.Range("B1049").Value = 'San Pablo Del Monte'
$ws.Range("B1054").Value = 'Tepetitla De Lardizábal'
$ws.Range("B1062").Value = 'Ziltlaltépec De Trinidad Sánchez Santos'
$ws.Range("B1069").Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range("B1073").Value = 'Amatlán De Los Reyes'
$ws.Range("B1080").Value = 'Boca Del Río'
$ws.Range("B1081").Value = 'Camarón De Tejeda'
$ws.Range("B1093").Value = 'Cosamaloapan De Carpio'
$ws.Range("B1094").Value = 'Cosautlán De Carvajal'
$ws.Range("B1106").Value = 'Hueyapan De Ocampo'
$ws.Range("B1107").Value = 'Huiloapan De Cuauhtémoc'
$ws.Range("B1108").Value = 'Ignacio De La Llave'
$ws.Range("B1111").Value = 'Ixhuacán De Los Reyes'
$ws.Range("B1112").Value = 'Ixhuatlán De Madero'
$ws.Range("B1113").Value = 'Ixhuatlán Del Café'
$ws.Range("B1122").Value = 'Juchique De Ferrer'
$ws.Range("B1126").Value = 'Lerdo De Tejada'
$ws.Range("B1128").Value = 'Martínez De La Torre'
$ws.Range("B1130").Value = 'Medellín De Bravo'
$ws.Range("B1141").Value = 'Ozuluama De Mascareñas'
$ws.Range("B1145").Value = 'Paso De Ovejas'
$ws.Range("B1146").Value = 'Paso Del Macho'
$ws.Range("B1149").Value = 'Poza Rica De Hidalgo'
$ws.Range("B1155").Value = 'Sayula De Alemán'
$ws.Range("B1156").Value = 'Soledad De Doblado'
$ws.Range("B1175").Value = 'Tlacotepec De Mejía'
$ws.Range("B1186").Value = 'Vega De Alatorre'
$ws.Range("B1219").Value = 'Moyahua De Estrada'
$ws.Range("B1220").Value = 'Nochistlán De Mejía'
$ws.Range("B1221").Value = 'Noria De Ángeles'
$ws.Range("B1228").Value = 'Teúl De González Ortega'
$ws.Range("B1229").Value = 'Tlaltenango De Sánchez Román'
$ws.Range("B1231").Value = 'Trinidad García De La Cadena'
$ws.Range("B1233").Value = 'Villa De Cos'

# --- Fix a floating point rounding artifact in D928 ---
$ws.Range("D928").Value = 0.009316064530788455

# --- Remove trailing metadata/footnote rows (1242-1246), keeping the Total row at 1240 ---
$ws.Rows("1242:1246").Delete()
